$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two mailto hyperlinks that used to live on B2/C2.
$ws.Hyperlinks.Delete()

# Swap the team names between row 2 and row 3 (Toronto Maple Leafs now goes
# with row 2's "leagueHomepageValidationTest"/"NHL Hockey" pairing, Montreal
# Canadiens moves down to row 3).
$c2 = $ws.Range("C2").Value2
$c3 = $ws.Range("C3").Value2
$ws.Range("C2").Value2 = $c3
$ws.Range("C3").Value2 = $c2

# Row 3 (B3/C3) already carries the "text, left/top aligned" formatting that
# used to be reserved for that row; copy it onto B2/C2 as well so every data
# row in column B/C is formatted identically.
$ws.Range("B3").Copy() | Out-Null
$ws.Range("B2").PasteSpecial(-4122) | Out-Null
$ws.Range("C3").Copy() | Out-Null
$ws.Range("C2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Drop the now-unused "addInfo" column entirely.
$ws.Columns.Item(4).Delete()
